$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B7: change from text "3" to numeric 3
$ws.Range("B7").Value = 3

# Add new row 8 data
$ws.Range("A8").Value = "Ying Tang"

# B8 stays as text "4" (not numeric)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "4"

$ws.Range("C8").Value = "As you suggested"
$ws.Range("D8").Value = "ACK"
$ws.Range("E8").Value = "EXP"
$ws.Range("F8").Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Range("G8").Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Range("H8").Value = "As you suggested, I did run comparison tests and I will present the results here."
